# Reorder the "Recorded By" (column G) comma-separated author lists:
# move the last item in the list to the front (rotate right by one),
# but only for cells that contain more than one comma-separated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ',\s*'
    if ($parts.Count -le 1) { continue }

    $lastIndex = $parts.Count - 1
    $newParts = @($parts[$lastIndex]) + $parts[0..($lastIndex - 1)]
    $newValue = [string]::Join(", ", $newParts)

    $cell.Value = $newValue
}
